$wb = $excel.ActiveWorkbook

$wsProperty = $wb.Worksheets.Item("Property")
$wsBroker = $wb.Worksheets.Item("Broker")

# Data change: Property sheet, Beds (E2) 5 -> 3
$wsProperty.Range("E2").Value = 3

# Update active selections to match the edited workbook state
$wsProperty.Activate()
$wsProperty.Range("E3").Select() | Out-Null

$wsBroker.Activate()
$wsBroker.Range("D29").Select() | Out-Null

$wsProperty.Activate()
